$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Consolidate the two tag rows (FinOps/Value04a and FinOps2/Value04b) into a
# single row: TagName becomes "FinOps3" with no TagValue, and the second
# data row is removed entirely.
$ws.Range("B2").Value = "FinOps3"
$ws.Range("C2").ClearContents()
$ws.Rows.Item(3).Delete()

# Move the active selection as recorded after the edits.
$ws.Range("C7").Select()
